$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new wishlist entry (row 44): The Complete Jacques Tati
# Set cells in an order that reproduces the shared-string insertion order
# observed in the target workbook (Link, Image, Name, Price).
$ws.Range("C44").Value = "https://www.criterion.com/boxsets/1069-the-complete-jacques-tati"
$ws.Range("B44").Value = "https://s3.amazonaws.com/criterion-production/product_images/1826-74abafd9cdfb81dc87d1738ec3c99693/VZ6wkU3rR5LF2srqIn918sRLIBPOl8_large.jpg"
$ws.Range("A44").Value = "The Complete Jacques Tati"
$ws.Range("D44").Value = "100 USD"

# Reflect the final selected cell as captured in the saved workbook
$ws.Range("E42").Select() | Out-Null
